$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Sheet1 -> Siddarth Kaul
$ws.Name = "Siddarth Kaul"

# Insert a new first column (A) -- everything (teamName..result) shifts
# right by one column, from A:L to B:M.
$ws.Columns.Item(1).Insert()

# Insert a new row 2 -- the existing (single) data row shifts down to
# row 3, freeing up row 2 for the newly-scraped match.
$ws.Rows.Item(2).Insert()

# ---- Header row ----
$ws.Range("A1").Value = "matchNo"

# ---- Row 2: new match (49th) ----
$ws.Range("A2").Value = "49th"
$ws.Range("B2").Value = "Sunrisers Hyderabad"
$ws.Range("C2").Value = "Siddarth Kaul"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'7"
$ws.Range("F2").Value = "'5"
$ws.Range("G2").Value = "'1"
$ws.Range("H2").Value = "'0"
$ws.Range("I2").Value = "'140.00"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Dubai (DSC)"
$ws.Range("L2").Value = "October 03"
$ws.Range("M2").Value = "KKR won by 6 wickets (with 2 balls remaining)"

# ---- Row 3: pre-existing match (55th) -- data (B3:M3) already shifted
# into place by the column/row inserts above; only matchNo is new. ----
$ws.Range("A3").Value = "55th"
